$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T1").Value = 0.99869362007748852
$ws.Range("AG1").Value = 0.77981146120731215
$ws.Range("AN1").Value = 0.93150182725711295
$ws.Range("BI2").Value = 0.62664343123870581
$ws.Range("AK3").Value = 0.98016668558573472
$ws.Range("BN3").Value = 0.73786472592550978
$ws.Range("F4").Value = 0.95039007782806939
$ws.Range("O4").Value = 0.92890721881526717
$ws.Range("AN5").Value = 0.94998839133393742
$ws.Range("BK5").Value = 0.8588220281714215
$ws.Range("BL5").Value = 0.62571103840130859
$ws.Range("E6").Value = 0.81575561939861851
$ws.Range("BM8").Value = 0.93288120786315054
$ws.Range("B9").Value = 0.72883842832974377
$ws.Range("K9").Value = 0.95663199835789992
$ws.Range("K10").Value = 0.74816072477451034
$ws.Range("Y10").Value = 0.70216878518288683
$ws.Range("AR10").Value = 0.92999353554773223
$ws.Range("AW10").Value = 0.78733227599403421
$ws.Range("P11").Value = 0.90414801123897792
$ws.Range("Q11").Value = 0.75019741595623879
$ws.Range("BD11").Value = 0.89599671328037456
$ws.Range("N12").Value = 0.88333120275532984
$ws.Range("X12").Value = 0.81472839143972553
$ws.Range("N13").Value = 0.94404371163443268
$ws.Range("Z13").Value = 0.85758688003321681
$ws.Range("AC14").Value = 0.98990153019225313
$ws.Range("AS14").Value = 0.83350285841163252
$ws.Range("C16").Value = 0.80605593920112573
$ws.Range("O16").Value = 0.62307939286471226
$ws.Range("AH16").Value = 0.89774102219073981
$ws.Range("BE16").Value = 0.8854236654860308
$ws.Range("BP16").Value = 0.81162623531102285
$ws.Range("W17").Value = 0.87141241412105952
$ws.Range("T18").Value = 0.61229010834648001
$ws.Range("AJ18").Value = 0.91139086023323568
$ws.Range("Q19").Value = 0.90730687479497751
$ws.Range("T19").Value = 0.99512131896038403
$ws.Range("U19").Value = 0.81854935398981898
$ws.Range("AR19").Value = 0.61109199562700334
$ws.Range("AV19").Value = 0.67410244795930718
$ws.Range("BB20").Value = 0.91427612589555074
$ws.Range("AB21").Value = 0.98844187592679833
$ws.Range("AJ21").Value = 0.69233154821932796
$ws.Range("O22").Value = 0.60794571985909829
$ws.Range("AA22").Value = 0.81675974679442986
$ws.Range("AI22").Value = 0.98869830784186163
$ws.Range("V23").Value = 0.79330454251787308
$ws.Range("BP23").Value = 0.91443981236203986
$ws.Range("D24").Value = 0.97973919275207499
$ws.Range("BE24").Value = 0.97251528209599525
$ws.Range("BG24").Value = 0.58710607765186218
$ws.Range("R25").Value = 0.83447362548418713
$ws.Range("AA26").Value = 0.75947555550222967
$ws.Range("AM26").Value = 0.79588381114698703
$ws.Range("AJ27").Value = 0.88615516104249537
$ws.Range("Z28").Value = 0.92929929357605057
$ws.Range("BD28").Value = 0.79087350698500147
$ws.Range("I29").Value = 0.98424540265866645
$ws.Range("AF29").Value = 0.89785474085875472
$ws.Range("AQ30").Value = 0.63560286805506694
$ws.Range("AB31").Value = 0.91249449191980814
$ws.Range("AS31").Value = 0.72614002866915439
$ws.Range("H32").Value = 0.84812734699671055
$ws.Range("AL32").Value = 0.86605839272339968
$ws.Range("AE33").Value = 0.95122432773898891
$ws.Range("AN33").Value = 0.83516690096664337
$ws.Range("B35").Value = 0.89781908666433641
$ws.Range("AH35").Value = 0.84324955055266981
$ws.Range("AN35").Value = 0.91181949466678347
$ws.Range("AX35").Value = 0.6014291781405936
$ws.Range("AD36").Value = 0.67549835441480144
$ws.Range("AE36").Value = 0.81794177123639122
$ws.Range("AL36").Value = 0.69803996706879501
$ws.Range("G37").Value = 0.96768511323831108
$ws.Range("I37").Value = 0.92001805703117934
$ws.Range("AZ37").Value = 0.95821136464703849
$ws.Range("BD37").Value = 0.79209202072541807
$ws.Range("AD38").Value = 0.85403232933752904
$ws.Range("AK38").Value = 0.65866437802670974
$ws.Range("AO39").Value = 0.93002842429591748
$ws.Range("AZ39").Value = 0.64070501939003377
$ws.Range("AB40").Value = 0.88399032378371223
$ws.Range("AL40").Value = 0.94711626225650236
$ws.Range("AI42").Value = 0.79013735957995723
$ws.Range("J43").Value = 0.6401599341766071
$ws.Range("AK43").Value = 0.94727587483511333
$ws.Range("AG44").Value = 0.97383374312249282
$ws.Range("L46").Value = 0.86236831054421681
$ws.Range("V47").Value = 0.72840251991051708
$ws.Range("AS47").Value = 0.89164552420792997
$ws.Range("BF47").Value = 0.91616520122545486
$ws.Range("X48").Value = 0.93088227900451281
$ws.Range("AT48").Value = 0.80923896155635955
$ws.Range("AK49").Value = 0.71035692903644221
$ws.Range("BM49").Value = 0.81338429077017516
$ws.Range("AW50").Value = 0.67820610169304119
$ws.Range("W51").Value = 0.85772698934845781
$ws.Range("BF51").Value = 0.9520690451849485
$ws.Range("H52").Value = 0.76839543797397813
$ws.Range("T52").Value = 0.59918107381604513
$ws.Range("AE52").Value = 0.99232817858182831
$ws.Range("AX52").Value = 0.67082199545025589
$ws.Range("Z53").Value = 0.61147474877228292
$ws.Range("BN53").Value = 0.61304710857809497
$ws.Range("AA54").Value = 0.65197757607803519
$ws.Range("BA54").Value = 0.70253255448120799
$ws.Range("BF54").Value = 0.72255446803517387
$ws.Range("BC56").Value = 0.67990899838365593
$ws.Range("BG57").Value = 0.95750912366456453
$ws.Range("N58").Value = 0.70540086010488956
$ws.Range("AI59").Value = 0.58736283627340546
$ws.Range("AP59").Value = 0.94016640628479209
$ws.Range("BM59").Value = 0.98670008466032011
$ws.Range("P60").Value = 0.90471576688424526
$ws.Range("Q60").Value = 0.85873877032464241
$ws.Range("AH60").Value = 0.97820461008365034
$ws.Range("AO60").Value = 0.99062033294133034
$ws.Range("AS60").Value = 0.95941178893828249
$ws.Range("BJ60").Value = 0.85357868949637949
$ws.Range("AF61").Value = 0.80808182861523126
$ws.Range("AK61").Value = 0.65387964140097532
$ws.Range("AO62").Value = 0.95980197421552438
$ws.Range("AT62").Value = 0.65108326979826869
$ws.Range("O63").Value = 0.80519427170120861
$ws.Range("BM63").Value = 0.95308687499101197
$ws.Range("F64").Value = 0.90667984740426766
$ws.Range("AU64").Value = 0.90860457997812061
$ws.Range("BJ64").Value = 0.99914791673159287
$ws.Range("G65").Value = 0.74959468578898558
$ws.Range("BN65").Value = 0.79774828956945409
$ws.Range("A67").Value = 0.84781125519940559
$ws.Range("F67").Value = 0.92165483125620118
$ws.Range("D68").Value = 0.8935977463630671
$ws.Range("BA68").Value = 0.993487217442107
$ws.Range("BC68").Value = 0.73255708316696877
